$wb = $excel.ActiveWorkbook

# Sheet "展览" — update "想去人数" (column F) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 940
$ws1.Range("F6").Value = 51

# Sheet "全部类型" — same events are duplicated here, update accordingly
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 940
$ws4.Range("F7").Value = 51
